$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-08-01 Thursday" "2024-08-02 Friday"

Replace-Text "23÷7=" "76÷7="
Replace-Text "40÷8=" "48÷3="
Replace-Text "97÷8=" "24÷4="
Replace-Text "82÷8=" "87÷5="
Replace-Text "98÷2=" "76÷6="

Replace-Text "91÷2=" "19÷5="
Replace-Text "18÷2=" "50÷6="
Replace-Text "69÷7=" "79÷8="
Replace-Text "85÷7=" "23÷6="
Replace-Text "67÷2=" "76÷4="

Replace-Text "37÷5=" "85÷2="
Replace-Text "92÷7=" "75÷6="
Replace-Text "32÷7=" "43÷9="
Replace-Text "86÷8=" "76÷6="
Replace-Text "17÷9=" "67÷6="

Replace-Text "83÷9=" "25÷7="
Replace-Text "12÷4=" "69÷2="
Replace-Text "43÷4=" "91÷6="
Replace-Text "52÷3=" "41÷8="
Replace-Text "94÷8=" "56÷2="

Replace-Text "84÷5=" "34÷5="
Replace-Text "81÷7=" "88÷9="
Replace-Text "13÷8=" "26÷5="
Replace-Text "75÷5=" "96÷9="
Replace-Text "50÷2=" "89÷6="
